$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two PRICE values that were filtered out as manual/duplicate entries ---
# (row 10 and row 30 had a PRICE in column C that gets removed)
$ws.Cells.Item(10, 3).ClearContents()
$ws.Cells.Item(30, 3).ClearContents()

# --- New comparison table in columns F:J, rows 40-45 ---
# Headers (row 40): F=IDEALISTA, G=REAL, H=PRED, I=Diff_Real, J=Diff_IDE
# Shared-string insertion order observed: PRED, REAL, IDEALISTA, Diff_Real, Diff_IDE
$ws.Cells.Item(40, 8).Value = "PRED"
$ws.Cells.Item(40, 7).Value = "REAL"
$ws.Cells.Item(40, 6).Value = "IDEALISTA"
$ws.Cells.Item(40, 9).Value = "Diff_Real"
$ws.Cells.Item(40, 10).Value = "Diff_IDE"

# Row 41 (single, non-shared formulas)
$ws.Cells.Item(41, 6).Value = 493000
$ws.Cells.Item(41, 7).Value = 485000
$ws.Cells.Item(41, 8).Value = 641511.9
$ws.Range("I41").Formula = "=(100*(G41/H41))"
$ws.Range("J41").Formula = "=(100*(F41/H41))"

# Rows 42-45 (data)
$ws.Cells.Item(42, 6).Value = 1665000
$ws.Cells.Item(42, 7).Value = 1500000
$ws.Cells.Item(42, 8).Value = 1750086.6

$ws.Cells.Item(43, 6).Value = 153000
$ws.Cells.Item(43, 7).Value = 149900
$ws.Cells.Item(43, 8).Value = 186867.78

$ws.Cells.Item(44, 6).Value = 738000
$ws.Cells.Item(44, 7).Value = 600000
$ws.Cells.Item(44, 8).Value = 743707.94

$ws.Cells.Item(45, 6).Value = 268000
$ws.Cells.Item(45, 7).Value = 250000
$ws.Cells.Item(45, 8).Value = 325048.16

# Rows 42-45 share one formula group each for columns I and J
$ws.Range("I42:I45").Formula = "=(100*(G42/H42))"
$ws.Range("J42:J45").Formula = "=(100*(F42/H42))"

# --- Column widths for the new table ---
$ws.Columns.Item(6).ColumnWidth = 9.92
$ws.Columns.Item(8).ColumnWidth = 8.92

# --- View state: selection + scroll position ---
$ws.Range("L46").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
